$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3487.087
$ws.Range("J17").Value = 3103.85
$ws.Range("L17").Value = 9311.549999999999
$ws.Range("N17").Value = -9647.549999999999
$ws.Range("H125").Value = 874.4286
$ws.Range("I125").Value = 650.5833
$ws.Range("J125").Value = 2217.5
$ws.Range("K125").Value = 5855.2497
$ws.Range("L125").Value = 19957.5
$ws.Range("M125").Value = -3395.2497
$ws.Range("N125").Value = -24877.5
$ws.Range("H132").Value = 1022.57776
$ws.Range("I132").Value = 934.9729599999999
$ws.Range("K132").Value = 2804.91888
$ws.Range("M132").Value = -274.9188799999997
$ws.Range("H137").Value = 27321.053
$ws.Range("I137").Value = 747.1905
$ws.Range("J137").Value = 60147.59
$ws.Range("K137").Value = 2241.5715
$ws.Range("L137").Value = 180442.77
$ws.Range("M137").Value = 308.4285
$ws.Range("N137").Value = -185542.77
$ws.Range("H138").Value = 2324.3774
$ws.Range("I138").Value = 2468.8572
$ws.Range("J138").Value = 2244.111
$ws.Range("K138").Value = 7406.571599999999
$ws.Range("L138").Value = 6732.333
$ws.Range("M138").Value = -2266.571599999999
$ws.Range("N138").Value = -17012.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2292.21
$ws.Range("I32").Value = 2060.6355
$ws.Range("J32").Value = 7850
$ws.Range("K32").Value = 2060.6355
$ws.Range("L32").Value = 7850
$ws.Range("M32").Value = -1773.6355
$ws.Range("N32").Value = -8424
$ws.Range("H61").Value = 38896.41
$ws.Range("I61").Value = 45628.945
$ws.Range("J61").Value = 8600
$ws.Range("K61").Value = 45628.945
$ws.Range("L61").Value = 8600
$ws.Range("M61").Value = -45416.945
$ws.Range("N61").Value = -9024
$ws.Range("H122").Value = 2485.4375
$ws.Range("I122").Value = 1365.7693
$ws.Range("J122").Value = 7337.3335
$ws.Range("K122").Value = 4097.3079
$ws.Range("L122").Value = 22012.0005
$ws.Range("M122").Value = -1647.3079
$ws.Range("N122").Value = -26912.0005
$ws.Range("H132").Value = 1873
$ws.Range("I132").Value = 1419.4584
$ws.Range("J132").Value = 2346.261
$ws.Range("K132").Value = 4258.3752
$ws.Range("L132").Value = 7038.782999999999
$ws.Range("M132").Value = -1728.3752
$ws.Range("N132").Value = -12098.783
$ws.Range("H136").Value = 38896.41
$ws.Range("I136").Value = 45628.945
$ws.Range("J136").Value = 8600
$ws.Range("K136").Value = 136886.835
$ws.Range("L136").Value = 25800
$ws.Range("M136").Value = -134336.835
$ws.Range("N136").Value = -30900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1709.2222
$ws.Range("I99").Value = 1672.875
$ws.Range("K99").Value = 1672.875
$ws.Range("M99").Value = -174.875
$ws.Range("H105").Value = 2568
$ws.Range("I105").Value = 2405.762
$ws.Range("K105").Value = 2405.762
$ws.Range("M105").Value = -658.7620000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2495.5833
$ws.Range("I31").Value = 1674.4286
$ws.Range("J31").Value = 3645.2
$ws.Range("K31").Value = 1674.4286
$ws.Range("L31").Value = 3645.2
$ws.Range("M31").Value = -1379.4286
$ws.Range("N31").Value = -4235.2
$ws.Range("H34").Value = 2495.5833
$ws.Range("I34").Value = 1674.4286
$ws.Range("J34").Value = 3645.2
$ws.Range("K34").Value = 1674.4286
$ws.Range("L34").Value = 3645.2
$ws.Range("M34").Value = -1472.4286
$ws.Range("N34").Value = -4049.2
$ws.Range("H51").Value = 30000
$ws.Range("H58").Value = 1360086
$ws.Range("I58").Value = 2175260
$ws.Range("J58").Value = 1462.4166
$ws.Range("K58").Value = 2175260
$ws.Range("L58").Value = 1462.4166
$ws.Range("M58").Value = -2175057
$ws.Range("N58").Value = -1868.4166
$ws.Range("H61").Value = 30000
$ws.Range("H94").Value = 1536.45
$ws.Range("I94").Value = 1490
$ws.Range("J94").Value = 1606.125
$ws.Range("K94").Value = 1490
$ws.Range("L94").Value = 1606.125
$ws.Range("M94").Value = -1039
$ws.Range("N94").Value = -2508.125
$ws.Range("H132").Value = 1523.6207
$ws.Range("I132").Value = 1141.35
$ws.Range("K132").Value = 3424.05
$ws.Range("M132").Value = -894.0499999999997
$ws.Range("H134").Value = 1326.7037
$ws.Range("I134").Value = 1175.2941
$ws.Range("J134").Value = 1584.1
$ws.Range("K134").Value = 3525.8823
$ws.Range("L134").Value = 4752.299999999999
$ws.Range("M134").Value = -990.8823000000002
$ws.Range("N134").Value = -9822.299999999999
$ws.Range("H135").Value = 35855.2
$ws.Range("J135").Value = 35855.2
$ws.Range("L135").Value = 35855.2
$ws.Range("N135").Value = -45995.2
$ws.Range("H136").Value = 1360086
$ws.Range("I136").Value = 2175260
$ws.Range("J136").Value = 1462.4166
$ws.Range("K136").Value = 6525780
$ws.Range("L136").Value = 4387.2498
$ws.Range("M136").Value = -6523230
$ws.Range("N136").Value = -9487.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 639
$ws.Range("J34").Value = 3000
$ws.Range("L34").Value = 9000
$ws.Range("N34").Value = -9168
$ws.Range("H122").Value = 1059.5385
$ws.Range("J122").Value = 1115.2609
$ws.Range("L122").Value = 10037.3481
$ws.Range("N122").Value = -14937.3481
$ws.Range("H131").Value = 26671.035
$ws.Range("I131").Value = 366
$ws.Range("J131").Value = 29827.64
$ws.Range("K131").Value = 1098
$ws.Range("L131").Value = 89482.92
$ws.Range("M131").Value = 3942
$ws.Range("N131").Value = -99562.92
$ws.Range("H137").Value = 5160.8696
$ws.Range("J137").Value = 7005
$ws.Range("L137").Value = 21015
$ws.Range("N137").Value = -31215

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2999
$ws.Range("I80").Value = 2999
$ws.Range("K80").Value = 2999
$ws.Range("M80").Value = -2001
$ws.Range("H83").Value = 2999
$ws.Range("I83").Value = 2999
$ws.Range("K83").Value = 14995
$ws.Range("M83").Value = -10003
$ws.Range("H102").Value = 2033.1428
$ws.Range("I102").Value = 2446.6
$ws.Range("K102").Value = 2446.6
$ws.Range("M102").Value = -824.5999999999999
$ws.Range("H122").Value = 3221.4443
$ws.Range("I122").Value = 3542.2856
$ws.Range("J122").Value = 2098.5
$ws.Range("K122").Value = 10626.8568
$ws.Range("L122").Value = 6295.5
$ws.Range("M122").Value = -8176.856800000001
$ws.Range("N122").Value = -11195.5
$ws.Range("H132").Value = 1541555.4
$ws.Range("I132").Value = 2264949
$ws.Range("J132").Value = 4344.25
$ws.Range("K132").Value = 6794847
$ws.Range("L132").Value = 13032.75
$ws.Range("M132").Value = -6792317
$ws.Range("N132").Value = -18092.75
$ws.Range("H136").Value = 8919.299999999999
$ws.Range("J136").Value = 8919.299999999999
$ws.Range("L136").Value = 26757.9
$ws.Range("N136").Value = -31857.9
$ws.Range("H141").Value = 44607.25
$ws.Range("J141").Value = 44607.25
$ws.Range("L141").Value = 44607.25
$ws.Range("N141").Value = -54967.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9624.714
$ws.Range("I16").Value = 13119.4
$ws.Range("K16").Value = 13119.4
$ws.Range("M16").Value = -12949.4
$ws.Range("H22").Value = 1224
$ws.Range("I22").Value = 833.8570999999999
$ws.Range("J22").Value = 1614.1428
$ws.Range("K22").Value = 833.8570999999999
$ws.Range("L22").Value = 1614.1428
$ws.Range("M22").Value = -538.8570999999999
$ws.Range("N22").Value = -2204.1428
$ws.Range("H27").Value = 1224
$ws.Range("I27").Value = 833.8570999999999
$ws.Range("J27").Value = 1614.1428
$ws.Range("K27").Value = 833.8570999999999
$ws.Range("L27").Value = 1614.1428
$ws.Range("M27").Value = -726.8570999999999
$ws.Range("N27").Value = -1828.1428
$ws.Range("H68").Value = 2374.875
$ws.Range("I68").Value = 1624.75
$ws.Range("K68").Value = 1624.75
$ws.Range("M68").Value = -875.75
$ws.Range("H71").Value = 2374.875
$ws.Range("I71").Value = 1624.75
$ws.Range("K71").Value = 8123.75
$ws.Range("M71").Value = -4379.75
$ws.Range("H82").Value = 1718.625
$ws.Range("I82").Value = 1718.625
$ws.Range("K82").Value = 1718.625
$ws.Range("M82").Value = -1357.625
$ws.Range("H85").Value = 1718.625
$ws.Range("I85").Value = 1718.625
$ws.Range("K85").Value = 1718.625
$ws.Range("M85").Value = -470.625
$ws.Range("H132").Value = 4430.3447
$ws.Range("I132").Value = 3724.818
$ws.Range("K132").Value = 11174.454
$ws.Range("M132").Value = -8644.454000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 29388
$ws.Range("J16").Value = 29388
$ws.Range("L16").Value = 29388
$ws.Range("N16").Value = -29972
$ws.Range("H96").Value = 10874
$ws.Range("I96").Value = 3300
$ws.Range("J96").Value = 11820.75
$ws.Range("K96").Value = 3300
$ws.Range("L96").Value = 11820.75
$ws.Range("M96").Value = -1927
$ws.Range("N96").Value = -14566.75
$ws.Range("H122").Value = 75044.09
$ws.Range("I122").Value = 91164.44500000001
$ws.Range("J122").Value = 2502.5
$ws.Range("K122").Value = 273493.335
$ws.Range("L122").Value = 7507.5
$ws.Range("M122").Value = -271043.335
$ws.Range("N122").Value = -12407.5
$ws.Range("H126").Value = 4387.0625
$ws.Range("J126").Value = 2387.6667
$ws.Range("L126").Value = 7163.000100000001
$ws.Range("N126").Value = -12103.0001
$ws.Range("H132").Value = 1356.766
$ws.Range("I132").Value = 1333.5667
$ws.Range("J132").Value = 1397.7059
$ws.Range("K132").Value = 4000.7001
$ws.Range("L132").Value = 4193.1177
$ws.Range("M132").Value = -1470.7001
$ws.Range("N132").Value = -9253.117699999999
